# [IMP] Adjust cheque register report
#
# - Report title (A2) translated from Thai to English.
# - A new "Lot Number" column is inserted between "Supplier Name" (D)
#   and the old "Cheque Number" column (now shifted from E to F).
# - Date columns (Posting Date, Cheque Received Date, Encashed Date,
#   Voided Date) get a DD/MM/YYYY number format, right aligned.
# - The Cheque Amount column gets an accounting-style numeric format,
#   right aligned.
# - General text columns become left aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before E ("Cheque Number" -> shifts to F, etc.)
$ws.Columns("E").Insert()

# 2. Translate the report title.
$ws.Range("A2").Value = "Cheque Register Report"

# 3. Fill in the header + width for the newly inserted "Lot Number" column.
$ws.Range("E10").Value = "Lot Number"
$ws.Columns("E").ColumnWidth = 25.42

# The inserted column duplicated the shaded "input" fill from D5 into E5;
# the new cell should stay unshaded/blank like a normal body cell.
$ws.Range("E5").Interior.ColorIndex = -4142
$ws.Range("E5").HorizontalAlignment = -4131

# 4. Number formats for the date columns (Posting Date, Cheque Received
#    Date, Encashed Date, Voided Date) and the Cheque Amount column.
$ws.Range("H1:H10").NumberFormat = "DD/MM/YYYY"
$ws.Range("H1:H10").HorizontalAlignment = -4152

$ws.Range("I1:I10").NumberFormat = "#,##0.00_);(#,##0.00)"
$ws.Range("I1:I10").HorizontalAlignment = -4152

$ws.Range("J1:J10").NumberFormat = "DD/MM/YYYY"
$ws.Range("J1:J10").HorizontalAlignment = -4152

$ws.Range("K1:K10").NumberFormat = "DD/MM/YYYY"
$ws.Range("K1:K10").HorizontalAlignment = -4152

$ws.Range("L1:L10").NumberFormat = "DD/MM/YYYY"
$ws.Range("L1:L10").HorizontalAlignment = -4152

# 5. The plain text columns (report header block + table text columns)
#    switch from "general" to explicit left alignment.
$ws.Columns("A").HorizontalAlignment = -4131
$ws.Columns("B").HorizontalAlignment = -4131
$ws.Columns("C").HorizontalAlignment = -4131
$ws.Columns("D").HorizontalAlignment = -4131
$ws.Columns("F").HorizontalAlignment = -4131
$ws.Columns("G").HorizontalAlignment = -4131
$ws.Columns("M").HorizontalAlignment = -4131

Write-Output "cheque register report layout updated"
